# Update cryptos list prices and 1h volume percentages
# (commit: "Updated cryptos list on Wed Nov 22 12:42:22 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.485.75"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "2.013.09"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.95"
$ws.Range("E5").Value = "  -10.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.599"
$ws.Range("E6").Value = "  -3.13%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.79"
$ws.Range("E8").Value = "  -3.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.369"
$ws.Range("E9").Value = "  -2.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.23"
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0747"
$ws.Range("E11").Value = "  -3.49%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "2.310.97"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.13"
$ws.Range("E15").Value = "  -8.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.762"
$ws.Range("E16").Value = "  -4.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.08"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").Value = "2.027.33"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "36.393.14"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("E20").Value = "  -3.88%  "
$ws.Range("D21").Value = "0.0₃0795"
$ws.Range("E21").Value = "  -5.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.36"
$ws.Range("E22").Value = "  +5.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "220.75"
$ws.Range("E23").Value = "  -5.33%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("E26").Value = "  -7.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.57"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.60"
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.129"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.37"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "18.88"
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.35"
$ws.Range("E33").Value = "  -5.39%  "
$ws.Range("E34").Value = "  -6.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.44"
$ws.Range("E35").Value = "  +3.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.24"
$ws.Range("E36").Value = "  -4.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.29"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.73"
$ws.Range("E40").Value = "  +3.19%  "
$ws.Range("E41").Value = "  -2.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0953"
$ws.Range("E42").Value = "  +2.90%  "
$ws.Range("D43").Value = "1.457.55"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.14"
$ws.Range("E44").Value = "  +37.02%  "
$ws.Range("E45").Value = "  -3.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.07"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("E47").Value = "  -7.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.40"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.92"
$ws.Range("E51").Value = "  -0.68%  "
